# Commit: "Deleted a row with lakshmi as email ID"
#
# The workbook has a login-attempt log on Sheet1:
#   Row1: username | password | res            (header)
#   Row2: salim@khan.com | salim1234 | Valid
#   Row3: lakshmi@yahoo.com | Lakkkk | Invalid   <- this row is removed
#   Row4: abc123@gmail.com | test@123 | Valid
#
# Deleting row 3 shifts old row 4 up to become the new row 3. That row's
# two cells (email + password) are both mailto hyperlinks, so the
# worksheet's Hyperlinks collection has to be rebuilt to drop the
# lakshmi link and re-target the ones that moved.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Remove the whole row (shifts rows 4.. up, fixes dimension/shared
#    strings/styles for the surviving cells automatically).
$ws.Rows(3).Delete()

# 2. The Hyperlinks collection does not follow the row shift by itself,
#    so it now points at stale/blank cells (A3 old-lakshmi-link, A4, B4).
#    Clear all of them and rebuild the three that should remain:
#      A2 -> salim@khan.com   (unchanged)
#      A3 -> abc123@gmail.com (was A4)
#      B3 -> test@123         (was B4)
$ws.Range("A2:B4").Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Cells.Item(2, 1), "mailto:salim@khan.com")
$ws.Hyperlinks.Add($ws.Cells.Item(3, 1), "mailto:abc123@gmail.com")
$ws.Hyperlinks.Add($ws.Cells.Item(3, 2), "mailto:test@123")

# 3. Re-adding a hyperlink resets the cell to a freshly-applied
#    "Hyperlink" style variant; nudging the underline back to the value
#    it already had collapses it onto the original shared cell style
#    instead of leaving a duplicate one behind.
$ws.Range("A2").Font.Underline = 2
$ws.Range("A3").Font.Underline = 2
$ws.Range("B3").Font.Underline = 2

# 4. Leave the selection where the author ended up.
$ws.Range("A13").Select()
